# Natmi following Dr Hou advice
# Replace the LR-pair results table (rows 2-17) on Sheet1 with the updated NATMI output:
# the sending/target cluster grid now spans all 4 clusters (ECs, FAPs, M2, sCs) instead of 3,
# and the underlying ligand/receptor/edge statistics were recomputed accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$allRows = @(
  @("ECs", "Dag1", "Lama4", "ECs", 3, 1, 11.36407333333333, 34.09222, 0.03922668827193482, 0.03922668827193482, 3, 1, 55.00581500000001, 165.017445, 0.2636796788967745, 0.2636796788967746, 625.0901154197667, 5625.8110387779, 0.01034328056772765, 0.01034328056772765),
  @("ECs", "Dag1", "Lama4", "FAPs", 3, 1, 11.36407333333333, 34.09222, 0.03922668827193482, 0.03922668827193482, 3, 1, 139.747495, 419.242485, 0.6699032567416481, 0.6699032567416482, 1588.100781329633, 14292.9070319667, 0.02627808622455855, 0.02627808622455855),
  @("ECs", "Dag1", "Lama4", "M2", 3, 1, 11.36407333333333, 34.09222, 0.03922668827193482, 0.03922668827193482, 1, 0.3333333333333333, 0.04249866666666666, 0.127496, 0.0002037245476720547, 0.0002037245476720547, 0.4829579645688888, 4.346621681119999, 0.000007991439324872614, 0.000007991439324872617),
  @("ECs", "Dag1", "Lama4", "sCs", 3, 1, 11.36407333333333, 34.09222, 0.03922668827193482, 0.03922668827193482, 3, 1, 13.81266366666667, 41.437991, 0.06621333981390531, 0.06621333981390533, 156.9681228366689, 1412.71310553002, 0.002597330040323755, 0.002597330040323755),
  @("FAPs", "Dag1", "Lama4", "ECs", 3, 1, 36.67092, 110.01276, 0.126581262307212, 0.126581262307212, 3, 1, 55.00581500000001, 165.017445, 0.2636796788967745, 0.2636796788967746, 2017.1138413998, 18154.0245725982, 0.03337690659951405, 0.03337690659951406),
  @("FAPs", "Dag1", "Lama4", "FAPs", 3, 1, 36.67092, 110.01276, 0.126581262307212, 0.126581262307212, 3, 1, 139.747495, 419.242485, 0.6699032567416481, 0.6699032567416482, 5124.669209345398, 46122.02288410859, 0.08479719986207017, 0.08479719986207018),
  @("FAPs", "Dag1", "Lama4", "M2", 3, 1, 36.67092, 110.01276, 0.126581262307212, 0.126581262307212, 1, 0.3333333333333333, 0.04249866666666666, 0.127496, 0.0002037245476720547, 0.0002037245476720547, 1.55846520544, 14.02618684896, 0.00002578771040729448, 0.00002578771040729448),
  @("FAPs", "Dag1", "Lama4", "sCs", 3, 1, 36.67092, 110.01276, 0.126581262307212, 0.126581262307212, 3, 1, 13.81266366666667, 41.437991, 0.06621333981390531, 0.06621333981390533, 506.5230843072399, 4558.707758765159, 0.008381368135220514, 0.008381368135220515),
  @("M2", "Dag1", "Lama4", "ECs", 3, 1, 4.688545666666666, 14.065637, 0.01618399617112621, 0.01618399617112621, 3, 1, 55.00581500000001, 165.017445, 0.2636796788967745, 0.2636796788967746, 257.8972755597183, 2321.075480037465, 0.004267390913669188, 0.004267390913669189),
  @("M2", "Dag1", "Lama4", "FAPs", 3, 1, 4.688545666666666, 14.065637, 0.01618399617112621, 0.01618399617112621, 3, 1, 139.747495, 419.242485, 0.6699032567416481, 0.6699032567416482, 655.2125121097715, 5896.912608987945, 0.01084171174213181, 0.01084171174213181),
  @("M2", "Dag1", "Lama4", "M2", 3, 1, 4.688545666666666, 14.065637, 0.01618399617112621, 0.01618399617112621, 1, 0.3333333333333333, 0.04249866666666666, 0.127496, 0.0002037245476720547, 0.0002037245476720547, 0.1992569394391111, 1.793312454952, 0.000003297077299488952, 0.000003297077299488953),
  @("M2", "Dag1", "Lama4", "sCs", 3, 1, 4.688545666666666, 14.065637, 0.01618399617112621, 0.01618399617112621, 3, 1, 13.81266366666667, 41.437991, 0.06621333981390531, 0.06621333981390533, 64.7613043794741, 582.8517394152669, 0.001071596438025722, 0.001071596438025723),
  @("sCs", "Dag1", "Lama4", "ECs", 3, 1, 236.9790546666667, 710.9371639999999, 0.8180080532497269, 0.8180080532497269, 3, 1, 55.00581500000001, 165.017445, 0.2636796788967745, 0.2636796788967746, 13035.22603986955, 117317.034358826, 0.2156921008158636, 0.2156921008158637),
  @("sCs", "Dag1", "Lama4", "FAPs", 3, 1, 236.9790546666667, 710.9371639999999, 0.8180080532497269, 0.8180080532497269, 3, 1, 139.747495, 419.242485, 0.6699032567416481, 0.6699032567416482, 33117.22925713472, 298055.0633142125, 0.5479862589128875, 0.5479862589128877),
  @("sCs", "Dag1", "Lama4", "M2", 3, 1, 236.9790546666667, 710.9371639999999, 0.8180080532497269, 0.8180080532497269, 1, 0.3333333333333333, 0.04249866666666666, 0.127496, 0.0002037245476720547, 0.0002037245476720547, 10.07129385126044, 90.64164466134399, 0.0001666483206403986, 0.0001666483206403987),
  @("sCs", "Dag1", "Lama4", "sCs", 3, 1, 236.9790546666667, 710.9371639999999, 0.8180080532497269, 0.8180080532497269, 3, 1, 13.81266366666667, 41.437991, 0.06621333981390531, 0.06621333981390533, 3273.31197815528, 29459.80780339752, 0.05416304520033532, 0.05416304520033533)
)

for ($i = 0; $i -lt $allRows.Count; $i++) {
    $r = $i + 2
    $rowData = $allRows[$i]
    $arr = New-Object 'object[,]' 1,20
    for ($j = 0; $j -lt 20; $j++) { $arr[0,$j] = $rowData[$j] }
    $ws.Range("A" + $r + ":T" + $r).Value = $arr
}
